# Insert a new weekly price record as row 90 in the "Bruselas (repollito)"
# sheet, pushing the former rows 90-93 down to 91-94 (dimension grows to
# A1:R94). This mirrors the weekly price-sheet update described by the
# commit message ("Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 90:93 down to 91:94, leaving a blank row 90 to populate.
$ws.Rows.Item(90).Insert()

$ws.Cells.Item(90, 1).Value  = 10
$ws.Cells.Item(90, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(90, 3).Value  = "La Araucanía"
$ws.Cells.Item(90, 4).Value  = 44753
$ws.Cells.Item(90, 5).Value  = 9
$ws.Cells.Item(90, 6).Value  = 100112035
$ws.Cells.Item(90, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(90, 8).Value  = "Sin especificar"
$ws.Cells.Item(90, 9).Value  = "Primera"
$ws.Cells.Item(90, 10).Value = 90
$ws.Cells.Item(90, 11).Value = 25000
$ws.Cells.Item(90, 12).Value = 26000
$ws.Cells.Item(90, 13).Value = 25444
$ws.Cells.Item(90, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(90, 15).Value = "Región Metropolitana"
$ws.Cells.Item(90, 16).Value = 2544
$ws.Cells.Item(90, 17).Value = 10
$ws.Cells.Item(90, 18).Value = "Hortaliza"
